$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the Predicted query text in B2: drop the appended
# "(rehabilitation OR ... AND ("Cervical Myelopathy"))" clause that was
# mistakenly concatenated onto the end of the query.
$newQuery = "`n" + '"spinal stenosis" OR "spinal surgery risks" OR "cervical myelopathy" OR "gait disturbance" OR "corticosteroid injections" OR "degenerative disc disease" OR "surgery for myelopathy" OR "myelopathy symptoms" OR "nerve root compression" OR "upper limb weakness" OR "mri cervical spine" OR "cervical radiculopathy" OR "cervical spondylotic myelopathy" OR "myelopathy assessment"' + "`n"
$ws.Range("B2").Value = $newQuery

# Setting a multi-line value can make the runtime auto-expand the row's
# height; re-running AutoFit restores the row to the (non-custom)
# default height so we don't leave a stray formatting change behind.
$ws.Rows.Item(2).AutoFit()

# --- Fix the Recall / Semantic Precision / Semantic F2 scores that had
# been swapped between the Predicted (row 2) and Baseline (row 3) rows.
$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2

$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2

$ws.Range("C2").Value = $c3
$ws.Range("D2").Value = $d3
$ws.Range("E2").Value = $e3

$ws.Range("C3").Value = $c2
$ws.Range("D3").Value = $d2
$ws.Range("E3").Value = $e2
